# Add surv bias rows for tuning purposes (commit: "Add surv bias for tuneing purposes")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("base")

# New rows 394-401: columns A..W hold literal values ("NULL" text maps to an
# existing shared string), columns X, Y, Z hold formulas identical in pattern
# to the rows already present in the sheet (e.g. row 393).
$data = @(
    ,(394,300,0,400,0.3,2,0.01,0.6,5,7,0.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T2")
    ,(395,300,0,400,0.99,2,0.01,0.6,5,7,0.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T2")
    ,(396,300,1,400,0.3,2,0.01,0.6,5,7,0.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T2")
    ,(397,300,1,400,0.99,2,0.01,0.6,5,7,0.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T2")
    ,(398,300,0,400,0.3,2,0.01,0.6,5,7,1.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T3")
    ,(399,300,0,400,0.99,2,0.01,0.6,5,7,1.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T3")
    ,(400,300,1,400,0.3,2,0.01,0.6,5,7,1.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T3")
    ,(401,300,1,400,0.99,2,0.01,0.6,5,7,1.5,0,30,0,"NULL",0,0,"NULL",0,0,0,0.2,0.2,"T3")
)

foreach ($rowVals in $data) {
    $r = $rowVals[0]

    $ws.Cells.Item($r, 1).Value  = $rowVals[1]   # A
    $ws.Cells.Item($r, 2).Value  = $rowVals[2]   # B
    $ws.Cells.Item($r, 3).Value  = $rowVals[3]   # C
    $ws.Cells.Item($r, 4).Value  = $rowVals[4]   # D
    $ws.Cells.Item($r, 5).Value  = $rowVals[5]   # E
    $ws.Cells.Item($r, 6).Value  = $rowVals[6]   # F
    $ws.Cells.Item($r, 7).Value  = $rowVals[7]   # G
    $ws.Cells.Item($r, 8).Value  = $rowVals[8]   # H
    $ws.Cells.Item($r, 9).Value  = $rowVals[9]   # I
    $ws.Cells.Item($r, 10).Value = $rowVals[10]  # J (survbias)
    $ws.Cells.Item($r, 11).Value = $rowVals[11]  # K
    $ws.Cells.Item($r, 12).Value = $rowVals[12]  # L
    $ws.Cells.Item($r, 13).Value = $rowVals[13]  # M
    $ws.Cells.Item($r, 14).Value = $rowVals[14]  # N
    $ws.Cells.Item($r, 15).Value = $rowVals[15]  # O
    $ws.Cells.Item($r, 16).Value = $rowVals[16]  # P
    $ws.Cells.Item($r, 17).Value = $rowVals[17]  # Q
    $ws.Cells.Item($r, 18).Value = $rowVals[18]  # R
    $ws.Cells.Item($r, 19).Value = $rowVals[19]  # S
    $ws.Cells.Item($r, 20).Value = $rowVals[20]  # T
    $ws.Cells.Item($r, 21).Value = $rowVals[21]  # U
    $ws.Cells.Item($r, 22).Value = $rowVals[22]  # V
    $ws.Cells.Item($r, 23).Value = $rowVals[23]  # W

    $ws.Cells.Item($r, 24).Formula = "=CONCATENATE(LOOKUP(D$r,info!`$C`$11:`$D`$19), F$r*100)"
    $ws.Cells.Item($r, 25).Formula = "=IF(AND(B$r=0,E$r=1),""F1"",IF(AND(B$r=0,E$r=2),""F2"",IF(AND(B$r=1,E$r=1),""M1"",IF(AND(B$r=1,E$r=2),""M2"",""?""))))"
    $ws.Cells.Item($r, 26).Formula = "=CONCATENATE(`$Y$r,""-"",`$W$r,""-"",`$X$r)"
}

# Select the newly added row (full row) as the active selection, matching
# the author's final cursor position after pasting the new rows in.
$ws.Range("A394:XFD394").Select()

# The "base" sheet becomes the active tab (tune was active before); this
# also flips tabSelected from "tune" to "base" and updates workbook.xml's
# bookViews/workbookView activeTab from 2 to 1.
$ws.Activate()
